$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Columns D, F, I, J, R, S (rows 1 header + 2-4 data) switch to Text number format
# and their data rows carry the numbers as text (shared-string) values instead
# of numeric literals.
$textCols = @("D", "F", "I", "J", "R", "S")
foreach ($col in $textCols) {
    $ws.Range($col + "1:" + $col + "4").NumberFormat = "@"
}

$ws.Range("D2").Value = "90001"
$ws.Range("D3").Value = "90001"
$ws.Range("D4").Value = "90001"

$ws.Range("F2").Value = "10011"
$ws.Range("F3").Value = "10011"
$ws.Range("F4").Value = "10011"

$ws.Range("I2").Value = "5"
$ws.Range("I3").Value = "5"
$ws.Range("I4").Value = "5"

$ws.Range("J2").Value = "5"
$ws.Range("J3").Value = "5"
$ws.Range("J4").Value = "5"

$ws.Range("R2").Value = "222"
$ws.Range("R3").Value = "222"
$ws.Range("R4").Value = "222"

$ws.Range("S2").Value = "999"
$ws.Range("S3").Value = "999"
$ws.Range("S4").Value = "999"

# V2 is replaced with a new order id (51480995) as part of the merge.
$ws.Range("V2").Value = "51480995"

# V4/W4's former contents ("51463320" / "05-16-2020") shift up onto row 3/4
# alignment stays the same (values unchanged on V3/V4/W3/W4), only the
# underlying style bookkeeping compacts - nothing else to do there.

# Selection moves to V2 as the last-edited cell.
$ws.Range("V2").Select()
